$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Once Color is picked, you can't change it" -- the three bug rows that
# were previously assigned to Steve (bugs #311, #289, #283, all about
# players being able to re-pick their color) are marked Done ("X")
# instead, same as the rest of the already-closed rows.
$ws.Range("C10").Value = "X"
$ws.Range("C12").Value = "X"
$ws.Range("C15").Value = "X"

# Scroll the view down a bit and move the active selection to B9, as it
# was left after reviewing/editing those rows.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B9").Select()
